# Change last transaction checkpoint
# - Memo (D2) changes from "payment" to "Test Payment"
# - toPayee (F2) changes from "electricity" to "water"
# - Selection on the "Global" sheet moves from A2 to G2
# - Column D width is narrowed slightly (content now shorter-ish / re-fit)
# - The originally active sheet ("Action1") must remain the active sheet

$wb = $excel.ActiveWorkbook
$wsGlobal = $wb.Worksheets.Item("Global")
$wsAction = $wb.Worksheets.Item("Action1")

# Update the transaction memo and payee values
$wsGlobal.Range("D2").Value = "Test Payment"
$wsGlobal.Range("F2").Value = "water"

# Re-fit column D to its new content width
$wsGlobal.Columns.Item(4).ColumnWidth = 11.35

# Move the checkpoint selection on the Global sheet to G2
$wsGlobal.Range("G2").Select()

# Restore the originally active sheet/tab
$wsAction.Activate()
